$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the two new columns -------------------------------------------------
# New column "Name" goes in before the existing "level" column (old C, new D).
$ws.Columns("C:C").Insert()
# New column "krId" goes in before the existing "KR phòng" column (old F, now G
# after the first insert shifted everything right by one).
$ws.Columns("G:G").Insert()

# --- Row 1 headers ---------------------------------------------------------------
$ws.Cells.Item(1, 2).Value2 = "employeeId"
$ws.Cells.Item(1, 3).Value2 = "Name"
$ws.Cells.Item(1, 4).Value2 = "level"
$ws.Cells.Item(1, 5).Value2 = "teamName"
$ws.Cells.Item(1, 6).Value2 = "Loại"
$ws.Cells.Item(1, 7).Value2 = "krId"
$ws.Cells.Item(1, 8).Value2 = "KR phòng"
$ws.Cells.Item(1, 9).Value2 = "KR team"
$ws.Cells.Item(1, 10).Value2 = "KR cá nhân"
$ws.Cells.Item(1, 11).Value2 = "Công thức tính"
$ws.Cells.Item(1, 12).Value2 = "Nguồn dữ liệu"
$ws.Cells.Item(1, 13).Value2 = "Định kỳ tính"
$ws.Cells.Item(1, 14).Value2 = "Đơn vị tính"
$ws.Cells.Item(1, 15).Value2 = "Điều kiện"
$ws.Cells.Item(1, 16).Value2 = "Norm"
$ws.Cells.Item(1, 17).Value2 = "% Trọng số chỉ tiêu"
$ws.Cells.Item(1, 18).Value2 = "Kết quả"
$ws.Cells.Item(1, 19).Value2 = "Tỷ lệ"
$ws.Cells.Item(1, 20).Value2 = "Tổng thời gian dự kiến/ ước tính công việc (giờ)"
$ws.Cells.Item(1, 21).Value2 = "Tổng thời gian thực hiện công việc thực tế (giờ)"
$ws.Cells.Item(1, 22).Value2 = "Note"

# --- Row 2 data --------------------------------------------------------------------
$ws.Cells.Item(2, 1).Value2 = 12
$ws.Cells.Item(2, 2).Value2 = 6
$ws.Cells.Item(2, 3).Value2 = "PHG"
$ws.Cells.Item(2, 4).Value2 = 1
$ws.Cells.Item(2, 5).Value2 = "python"
$ws.Cells.Item(2, 6).Value2 = "KPI"
$ws.Cells.Item(2, 7).Value2 = 4
$ws.Cells.Item(2, 8).Value2 = "tét2321"
$ws.Cells.Item(2, 9).Value2 = "tét2321"
$ws.Cells.Item(2, 10).Value2 = "tét2321"
$ws.Cells.Item(2, 11).Value2 = "Báo cáo được CBQL confirm"
$ws.Cells.Item(2, 12).Value2 = "email"
$ws.Cells.Item(2, 13).Value2 = "Tháng"
$ws.Cells.Item(2, 14).Value2 = "%"

# O2 must hold the literal text "=" rather than be parsed as a formula, so we
# use a leading apostrophe (quote-prefix) to force text entry, then reset the
# cell style back to Normal so the quote-prefix flag doesn't linger as an
# extra style.
$o2 = $ws.Cells.Item(2, 15)
$o2.Value2 = "'="
$o2.Style = "Normal"

$ws.Cells.Item(2, 16).Value2 = 100
$ws.Cells.Item(2, 17).Value2 = 78
$ws.Cells.Item(2, 18).Value2 = 0
$ws.Cells.Item(2, 19).Value2 = 0

# T2/U2 hold the digit-only text "168" (not a number), so again use the
# quote-prefix trick to keep it text, then clear the resulting style flag.
$t2 = $ws.Cells.Item(2, 20)
$t2.Value2 = "'168"
$t2.Style = "Normal"

$u2 = $ws.Cells.Item(2, 21)
$u2.Value2 = "'168"
$u2.Style = "Normal"

$ws.Cells.Item(2, 22).Value2 = "dfqwq"
